$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the blank spacer row (old row 3) so rows shift up.
$ws.Rows.Item(3).Delete()

# 2. Clear the subtitle text in (now) row 2 - entire row becomes blank.
$ws.Rows.Item(2).Clear()

# 2b. The title row no longer carries a styled "filler" cell next to it.
$ws.Range("B1").Clear()

# 3. The table is being trimmed down to a single (2014) data column. Column D
#    held that figure and also carried the "closing" border of the box, so
#    pull its formatting into column B before the value itself is moved.
$ws.Range("D4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("B4").Value2 = $ws.Range("D4").Value2

# 4. Drop the now unused columns C and D (delete right-to-left so indices
#    stay valid).
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(3).Delete()

# 5. Set the uniform row height (20.1pt) used throughout the refreshed layout.
$ws.Range("A1:A5").EntireRow.RowHeight = 20.1

# 6. Append two blank rows at the bottom of the table, matching the new layout.
$ws.Rows.Item(6).RowHeight = 20.1
$ws.Rows.Item(7).RowHeight = 20.1
